# Update the "Tools:" bullet on the Data slide to add SQLAlchemy to the
# existing Python/Pandas/Numpy tooling line.
$p = $ppt.ActivePresentation

$targetSlide = $null
$targetShape = $null
$targetParaIndex = -1

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            $paraCount = $tr.Paragraphs().Count
            for ($k = 1; $k -le $paraCount; $k++) {
                $para = $tr.Paragraphs($k, 1)
                $paraText = $para.Text.TrimEnd("`r")
                if ($paraText -eq "Python, Pandas, Numpy") {
                    $targetSlide = $s
                    $targetShape = $sh
                    $targetParaIndex = $k
                }
            }
        }
    }
}

if ($targetParaIndex -ge 1) {
    $tr = $targetShape.TextFrame.TextRange
    $para = $tr.Paragraphs($targetParaIndex, 1)
    $para.Text = "Python, Pandas, Numpy, SQLAlchemy"
}
